# Update cryptos list sheet with the latest scraped values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.581.51"
$ws.Range("E2").Value = "  +2.53%  "

# Row 3
$ws.Range("D3").Value = "2.698.77"
$ws.Range("E3").Value = "  +2.31%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.577"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.86%  "

# Row 9
$ws.Range("D9").Value = "2.719.88"
$ws.Range("E9").Value = "  +2.48%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.106"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.44%  "

# Row 12
$ws.Range("E12").Value = "  +0.34%  "

# Row 13
$ws.Range("E13").Value = "  +3.03%  "

# Row 14
$ws.Range("D14").Value = "3.176.58"
$ws.Range("E14").Value = "  +2.46%  "

# Row 15
$ws.Range("D15").Value = "60.565.12"
$ws.Range("E15").Value = "  +2.51%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.02%  "

# Row 17
$ws.Range("D17").Value = "2.716.88"
$ws.Range("E17").Value = "  +2.86%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000138"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.12%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "346.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.35%  "

# Row 20
$ws.Range("E20").Value = "  +0.04%  "

# Row 21
$ws.Range("E21").Value = "  +2.88%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.42%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.57%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.422"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.70%  "

# Row 26
$ws.Range("E26").Value = "  +2.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.993"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0819"
$ws.Range("E28").Value = "  +1.34%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.00%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.00%  "

# Row 32
$ws.Range("E32").Value = "  +0.89%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.18%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "

# Row 35
$ws.Range("E35").Value = "  +6.14%  "

# Row 36
$ws.Range("E36").Value = "  +8.50%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.939"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.50%  "

# Row 38 -> Fetch.AI
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.874"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.99%  "

# Row 39 -> Stacks
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.67%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.91%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.53%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "282.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.19%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.17%  "

# Row 44 -> Maker
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.145.83"
$ws.Range("E44").Value = "  +7.80%  "

# Row 45 -> Stellar
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0987"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "

# Row 46 -> FirstDigitalUSD
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.995"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "

# Row 47
$ws.Range("E47").Value = "  -0.34%  "

# Row 48
$ws.Range("E48").Value = "  +2.19%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.13%  "

# Row 50
$ws.Range("E50").Value = "  +1.85%  "

# Row 51
$ws.Range("E51").Value = "  +0.90%  "
